# Add a new "Table of Contents" tab as the first sheet in the workbook,
# describing the various data tables (sheets / source files).

$wb = $excel.ActiveWorkbook

# Create the new sheet by copying an existing sheet (caper_business) so that
# it inherits the workbook's normal sheet formatting (default row height,
# x14ac namespace, etc.) without carrying over any page-setup settings.
# The copy is placed before the first sheet, so it becomes the new sheet #1.
$wb.Worksheets.Item(2).Copy($wb.Worksheets.Item(1))
$toc = $wb.Worksheets.Item(1)
$toc.Name = "Table of Contents"

# Wipe out all of the copied rows/content - we only want the formatting
# environment (sheetFormatPr, namespaces) from the source sheet.
$toc.UsedRange.EntireRow.Delete() | Out-Null

# Column widths (characters), matching the authored table layout.
$toc.Columns.Item(1).ColumnWidth = 12.76
$toc.Columns.Item(2).ColumnWidth = 46.76
$toc.Columns.Item(3).ColumnWidth = 28.1
$toc.Columns.Item(4).ColumnWidth = 17.93
$toc.Columns.Item(5).ColumnWidth = 66.76

# Header row
$toc.Range("A1").Value = "File"
$toc.Range("B1").Value = "Description"
$toc.Range("C1").Value = "Unit of Analysis"
$toc.Range("D1").Value = "Merge Variable(s)"
$toc.Range("E1").Value = "Multiple Files"

# Row 3 - army_master
$toc.Range("A3").Value = "army_master"
$toc.Range("B3").Value = " Personnnel Files"
$toc.Range("C3").Value = "individual servicemember"
$toc.Range("D3").Value = "PID_PDE SNPSHT_DT"
$toc.Range("E3").Value = "Separate Files for Army, and other services"

# Row 4 - PID_NPI
$toc.Range("A4").Value = "PID_NPI"
$toc.Range("B4").Value = "crosswalk between military ID and NPI for anyone (military or civilian) with a patient encouner since third quarter 2012"
$toc.Range("C4").Value = "Individual  Provider"
$toc.Range("D4").Value = "PID_PDE OR PROVNPI"
$toc.Rows.Item(4).RowHeight = 43.5
$toc.Range("B4").WrapText = $true

# Row 5 - caper_patient
$toc.Range("A5").Value = "caper_patient"
$toc.Range("B5").Value = "Comprehensive Ambulatory Provider Encounter Record"
$toc.Range("C5").Value = "provider encounter (multiple cpt)"
$toc.Range("D5").Value = "encounter_key"
$toc.Range("E5").Value = "4 Files each for Army, Non-Army, and Dependents due to STATA size limitations"

# Row 6 - caper_business
$toc.Range("A6").Value = "caper_business"
$toc.Range("B6").Value = "Information on RVU's and appointment info for provider level encounters"
$toc.Range("C6").Value = "provider encounter"
$toc.Range("D6").Value = "encounter_key"
$toc.Range("E6").Value = "4 Files each for Army, Non-Army, and Dependents due to STATA size limitations"
$toc.Rows.Item(6).RowHeight = 29
$toc.Range("B6").WrapText = $true

# Row 7 - sidr_patient
$toc.Range("A7").Value = "sidr_patient"
$toc.Range("B7").Value = "Standard Inpatient Data Record"
$toc.Range("C7").Value = "Inpatient Admission"
$toc.Range("D7").Value = "DMIS_PATIENT"
$toc.Range("E7").Value = "Separate Files for Army, Non-Army and Dependents"

# Row 8 - sidr_provider
$toc.Range("A8").Value = "sidr_provider"
$toc.Range("B8").Value = "provider information related to an inpatient stay"
$toc.Range("C8").Value = "Inpatient Admission"
$toc.Range("D8").Value = "DMIS_PATIENT"
$toc.Range("E8").Value = "Separate Files for Army, Non-Army and Dependents"
$toc.Range("B8").WrapText = $true

# Row 9 - TED-NI
$toc.Range("A9").Value = "TED-NI"
$toc.Range("B9").Value = "TRICARE Encounter Data; each record is a non-denied line item claim for services (other than hospital or institutional care for inpatients) provided in the private sector."
$toc.Range("C9").Value = "Provider claim (one cpt per claim)"
$toc.Range("D9").Value = "TEDNO, ADMTEDNO"
$toc.Range("E9").Value = "4 Files each for Army, Non-Army, and Dependents due to STATA size limitations"
$toc.Rows.Item(9).RowHeight = 58
$toc.Range("B9").WrapText = $true

# Row 10 - ted-I
$toc.Range("A10").Value = "ted-I"
$toc.Range("B10").Value = "TRICARE Encounter Data; each record is a hospital claim for services provided in the private sector. "
$toc.Range("C10").Value = "Inpatient Admission"
$toc.Range("D10").Value = "TEDNO, ADMTEDNO"
$toc.Range("E10").Value = "Separate Files for Army, Non-Army and Dependents"
$toc.Rows.Item(10).RowHeight = 29
$toc.Range("B10").WrapText = $true

# Row 11 - PDTS
$toc.Range("A11").Value = "PDTS"
$toc.Range("B11").Value = "Pharmacy Transactions - not complete before 2011"
$toc.Range("C11").Value = "Prescription Dispensed"
$toc.Range("E11").Value = "Separate Files for Army, Non-Army and Dependents"

# Make the Table of Contents the active sheet / selection, matching the
# authored view state.
$toc.Range("D17").Select() | Out-Null
